$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price + 1h volume change) per row.
# NumberFormat "@" (Text) is set before assigning the value so that
# numeric-looking strings (e.g. "0.999", "19.01") are stored as literal
# text - matching the workbook's inlineStr/shared-string cell type -
# instead of being auto-coerced into Excel numbers. Style is reset back
# to "Normal" right after so no stray number-format style is left on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "70.186.68"
Set-TextValue "E2" "  +0.69%  "
Set-TextValue "D3" "3.561.82"
Set-TextValue "E3" "  +0.86%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.18%  "
Set-TextValue "D5" "610.33"
Set-TextValue "E5" "  +4.39%  "
Set-TextValue "D6" "186.11"
Set-TextValue "E6" "  +1.48%  "
Set-TextValue "D7" "3.555.29"
Set-TextValue "E7" "  +0.95%  "
Set-TextValue "D8" "0.617"
Set-TextValue "E8" "  +0.91%  "
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "D10" "0.214"
Set-TextValue "E10" "  +9.14%  "
Set-TextValue "D11" "0.646"
Set-TextValue "E11" "  +0.61%  "
Set-TextValue "D12" "53.97"
Set-TextValue "E12" "  -0.03%  "
Set-TextValue "D13" "0.0000310"
Set-TextValue "E13" "  +2.07%  "
Set-TextValue "D14" "9.58"
Set-TextValue "E14" "  +1.57%  "
Set-TextValue "D15" "4.121.88"
Set-TextValue "E15" "  +0.68%  "
Set-TextValue "D16" "70.276.16"
Set-TextValue "E16" "  +0.84%  "
Set-TextValue "D17" "3.572.26"
Set-TextValue "E17" "  +1.37%  "
Set-TextValue "D18" "19.01"
Set-TextValue "E18" "  -1.30%  "
Set-TextValue "D19" "12.69"
Set-TextValue "E19" "  +2.81%  "
Set-TextValue "D20" "578.47"
Set-TextValue "E20" "  +7.90%  "
Set-TextValue "E21" "  +0.70%  "
Set-TextValue "D22" "0.994"
Set-TextValue "E22" "  -1.38%  "
Set-TextValue "D23" "17.42"
Set-TextValue "E23" "  -1.89%  "
Set-TextValue "E24" "  +4.17%  "
Set-TextValue "E25" "  +1.05%  "
Set-TextValue "D26" "94.22"
Set-TextValue "E26" "  -1.55%  "
Set-TextValue "D27" "2.95"
Set-TextValue "E27" "  -0.89%  "
Set-TextValue "D28" "10.95"
Set-TextValue "E28" "  -1.83%  "
Set-TextValue "D29" "9.39"
Set-TextValue "E29" "  +3.73%  "
Set-TextValue "D30" "32.28"
Set-TextValue "E30" "  +0.76%  "
Set-TextValue "D31" "7.05"
Set-TextValue "E31" "  -2.85%  "
Set-TextValue "D32" "12.25"
Set-TextValue "E32" "  -0.85%  "
Set-TextValue "D33" "0.115"
Set-TextValue "E33" "  +2.43%  "
Set-TextValue "D34" "64.01"
Set-TextValue "E34" "  -0.37%  "
Set-TextValue "D35" "3.72"
Set-TextValue "E35" "  +21.52%  "
Set-TextValue "D36" "3.19"
Set-TextValue "E36" "  +0.37%  "
Set-TextValue "B37" "Bittensor"
Set-TextValue "C37" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D37" "526.60"
Set-TextValue "E37" "  -3.31%  "
Set-TextValue "B38" "TheGraph"
Set-TextValue "C38" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D38" "0.406"
Set-TextValue "E38" "  -1.02%  "
Set-TextValue "E39" "  +0.21%  "
Set-TextValue "D40" "37.43"
Set-TextValue "E40" "  -1.49%  "
Set-TextValue "D41" "3.565.55"
Set-TextValue "E41" "  +6.32%  "
Set-TextValue "D42" "0.0₃0782"
Set-TextValue "E42" "  +3.33%  "
Set-TextValue "D43" "3.54"
Set-TextValue "E43" "  +5.09%  "
Set-TextValue "E44" "  +2.14%  "
Set-TextValue "D45" "0.0459"
Set-TextValue "E45" "  +4.89%  "
Set-TextValue "B46" "ThetaToken"
Set-TextValue "C46" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D46" "2.93"
Set-TextValue "E46" "  -0.64%  "
Set-TextValue "B47" "ApeXProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "3.44"
Set-TextValue "E47" "  -3.49%  "
Set-TextValue "D48" "0.140"
Set-TextValue "E48" "  +3.74%  "
Set-TextValue "D49" "9.18"
Set-TextValue "E49" "  +0.72%  "
Set-TextValue "E50" "  +0.22%  "
Set-TextValue "D51" "135.59"
Set-TextValue "E51" "  -1.07%  "
